$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 360.9  # was 557.62
$ws.Range("I17").Value = 213.33333  # was 176.66667
$ws.Range("J17").Value = 365.46393  # was 581.93616
$ws.Range("K17").Value = 639.99999  # was 530.00001
$ws.Range("L17").Value = 1096.39179  # was 1745.80848
$ws.Range("M17").Value = -471.99999  # was -362.00001
$ws.Range("N17").Value = -1432.39179  # was -2081.80848

# Row 64
$ws.Range("H64").Value = 128250  # was 38907.43
$ws.Range("I64").Value = 4333.3335  # was 3800.8
$ws.Range("J64").Value = 500000  # was 58411.11
$ws.Range("K64").Value = 4333.3335  # was 3800.8
$ws.Range("L64").Value = 500000  # was 58411.11
$ws.Range("M64").Value = -4085.3335  # was -3552.8
$ws.Range("N64").Value = -500496  # was -58907.11

# Row 67
$ws.Range("H67").Value = 128250  # was 38907.43
$ws.Range("I67").Value = 4333.3335  # was 3800.8
$ws.Range("J67").Value = 500000  # was 58411.11
$ws.Range("K67").Value = 4333.3335  # was 3800.8
$ws.Range("L67").Value = 500000  # was 58411.11
$ws.Range("M67").Value = -3475.3335  # was -2942.8
$ws.Range("N67").Value = -501716  # was -60127.11

# Row 132
$ws.Range("H132").Value = 174413.47  # was 202294.27
$ws.Range("I132").Value = 177455.81  # was 202294.27
$ws.Range("J132").Value = 1000  # was 0
$ws.Range("K132").Value = 532367.4299999999  # was 606882.8099999999
$ws.Range("L132").Value = 3000  # was 0
$ws.Range("M132").Value = -529837.4299999999  # was -604352.8099999999
$ws.Range("N132").Value = -8060  # new cell (previously empty)

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7228.0146  # was 6584.16
$ws.Range("I32").Value = 4615.4106  # was 4093.0952
$ws.Range("J32").Value = 18482.309  # was 19662.25
$ws.Range("K32").Value = 4615.4106  # was 4093.0952
$ws.Range("L32").Value = 18482.309  # was 19662.25
$ws.Range("M32").Value = -4328.4106  # was -3806.0952
$ws.Range("N32").Value = -19056.309  # was -20236.25

# Row 61
$ws.Range("H61").Value = 12822243  # was 9260859
$ws.Range("I61").Value = 17545546  # was 11495782
$ws.Range("J61").Value = 1845.7142  # was 1891.1428
$ws.Range("K61").Value = 17545546  # was 11495782
$ws.Range("L61").Value = 1845.7142  # was 1891.1428
$ws.Range("M61").Value = -17545334  # was -11495570
$ws.Range("N61").Value = -2269.7142  # was -2315.1428

# Row 74
$ws.Range("H74").Value = 6671909.5  # was 8702648
$ws.Range("I74").Value = 9091831  # was 15385616
$ws.Range("J74").Value = 17123.625  # was 14789.9
$ws.Range("K74").Value = 9091831  # was 15385616
$ws.Range("L74").Value = 17123.625  # was 14789.9
$ws.Range("M74").Value = -9090957  # was -15384742
$ws.Range("N74").Value = -18871.625  # was -16537.9

# Row 77
$ws.Range("H77").Value = 6671909.5  # was 8702648
$ws.Range("I77").Value = 9091831  # was 15385616
$ws.Range("J77").Value = 17123.625  # was 14789.9
$ws.Range("K77").Value = 45459155  # was 76928080
$ws.Range("L77").Value = 85618.125  # was 73949.5
$ws.Range("M77").Value = -45454787  # was -76923712
$ws.Range("N77").Value = -94354.125  # was -82685.5

# Row 122
$ws.Range("H122").Value = 1217.2903  # was 1190.5454
$ws.Range("I122").Value = 1268.6666  # was 1177.8334
$ws.Range("J122").Value = 1146.1538  # was 1224.4445
$ws.Range("K122").Value = 3805.9998  # was 3533.5002
$ws.Range("L122").Value = 3438.4614  # was 3673.3335
$ws.Range("M122").Value = -1355.9998  # was -1083.5002
$ws.Range("N122").Value = -8338.4614  # was -8573.333500000001

# Row 136
$ws.Range("H136").Value = 12822243  # was 9260859
$ws.Range("I136").Value = 17545546  # was 11495782
$ws.Range("J136").Value = 1845.7142  # was 1891.1428
$ws.Range("K136").Value = 52636638  # was 34487346
$ws.Range("L136").Value = 5537.142599999999  # was 5673.428400000001
$ws.Range("M136").Value = -52634088  # was -34484796
$ws.Range("N136").Value = -10637.1426  # was -10773.4284

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1756.5593  # was 1855.5555
$ws.Range("I86").Value = 1628.3658  # was 1750.1714
$ws.Range("J86").Value = 2048.5557  # was 2049.6843
$ws.Range("K86").Value = 1628.3658  # was 1750.1714
$ws.Range("L86").Value = 2048.5557  # was 2049.6843
$ws.Range("M86").Value = -505.3658  # was -627.1713999999999
$ws.Range("N86").Value = -4294.5557  # was -4295.6843

# Row 89
$ws.Range("H89").Value = 1756.5593  # was 1855.5555
$ws.Range("I89").Value = 1628.3658  # was 1750.1714
$ws.Range("J89").Value = 2048.5557  # was 2049.6843
$ws.Range("K89").Value = 8141.829  # was 8750.857
$ws.Range("L89").Value = 10242.7785  # was 10248.4215
$ws.Range("M89").Value = -2525.829  # was -3134.857
$ws.Range("N89").Value = -21474.7785  # was -21480.4215

# Row 134
$ws.Range("H134").Value = 51258.363  # was 39190.758
$ws.Range("I134").Value = 69973.125  # was 51233.773
$ws.Range("J134").Value = 1352.3334  # was 1341.2858
$ws.Range("K134").Value = 209919.375  # was 153701.319
$ws.Range("L134").Value = 4057.0002  # was 4023.8574
$ws.Range("M134").Value = -207384.375  # was -151166.319
$ws.Range("N134").Value = -9127.0002  # was -9093.857400000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1806.9688  # was 1701.129
$ws.Range("I31").Value = 1557.2  # was 1548.4
$ws.Range("J31").Value = 1920.5  # was 1773.8572
$ws.Range("K31").Value = 1557.2  # was 1548.4
$ws.Range("L31").Value = 1920.5  # was 1773.8572
$ws.Range("M31").Value = -1262.2  # was -1253.4
$ws.Range("N31").Value = -2510.5  # was -2363.8572

# Row 34
$ws.Range("H34").Value = 1806.9688  # was 1701.129
$ws.Range("I34").Value = 1557.2  # was 1548.4
$ws.Range("J34").Value = 1920.5  # was 1773.8572
$ws.Range("K34").Value = 1557.2  # was 1548.4
$ws.Range("L34").Value = 1920.5  # was 1773.8572
$ws.Range("M34").Value = -1355.2  # was -1346.4
$ws.Range("N34").Value = -2324.5  # was -2177.8572

# Row 62
$ws.Range("H62").Value = 2850  # was 2660
$ws.Range("I62").Value = 0  # was 2670
$ws.Range("J62").Value = 2850  # was 2650
$ws.Range("K62").Value = 0  # was 2670
$ws.Range("L62").Value = 2850  # was 2650
$ws.Range("M62").ClearContents()  # was -2046
$ws.Range("N62").Value = -4098  # was -3898

# Row 65
$ws.Range("H65").Value = 2850  # was 2660
$ws.Range("I65").Value = 0  # was 2670
$ws.Range("J65").Value = 2850  # was 2650
$ws.Range("K65").Value = 0  # was 13350
$ws.Range("L65").Value = 14250  # was 13250
$ws.Range("M65").ClearContents()  # was -10230
$ws.Range("N65").Value = -20490  # was -19490

# Row 132
$ws.Range("H132").Value = 1514.772  # was 1348.8055
$ws.Range("I132").Value = 1414.6  # was 1277.8387
$ws.Range("J132").Value = 2230.2856  # was 1788.8
$ws.Range("K132").Value = 4243.799999999999  # was 3833.5161
$ws.Range("L132").Value = 6690.8568  # was 5366.4
$ws.Range("M132").Value = -1713.799999999999  # was -1303.5161
$ws.Range("N132").Value = -11750.8568  # was -10426.4

$ws = $wb.Worksheets.Item("CUL")
# Row 105
$ws.Range("H105").Value = 302668670  # was 454001500
$ws.Range("J105").Value = 302668670  # was 454001500
$ws.Range("L105").Value = 908006010  # was 1362004500
$ws.Range("N105").Value = -908011252  # was -1362009742

# Row 110
$ws.Range("H110").Value = 2587.375  # was 2644.3333
$ws.Range("I110").Value = 2399.8572  # was 2487.375
$ws.Range("K110").Value = 7199.571599999999  # was 7462.125
$ws.Range("M110").Value = -3109.571599999999  # was -3372.125

# Row 131
$ws.Range("H131").Value = 2747.5576  # was 2736.6538
$ws.Range("I131").Value = 17100  # was 12975
$ws.Range("J131").Value = 1868.8368  # was 1883.4584
$ws.Range("K131").Value = 51300  # was 38925
$ws.Range("L131").Value = 5606.5104  # was 5650.3752
$ws.Range("M131").Value = -46260  # was -33885
$ws.Range("N131").Value = -15686.5104  # was -15730.3752

# Row 137
$ws.Range("H137").Value = 25964.111  # was 28383.781
$ws.Range("I137").Value = 2676.1538  # was 3119
$ws.Range("J137").Value = 35424.844  # was 36533.71
$ws.Range("K137").Value = 8028.4614  # was 9357
$ws.Range("L137").Value = 106274.532  # was 109601.13
$ws.Range("M137").Value = -2928.4614  # was -4257
$ws.Range("N137").Value = -116474.532  # was -119801.13

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 114.117645  # was 84.13043
$ws.Range("I2").Value = 54  # was 38.454544
$ws.Range("J2").Value = 181.75  # was 126
$ws.Range("K2").Value = 54  # was 38.454544
$ws.Range("L2").Value = 181.75  # was 126
$ws.Range("M2").Value = 59  # was 74.545456
$ws.Range("N2").Value = -407.75  # was -352

# Row 5
$ws.Range("H5").Value = 2090.5444  # was 2071.5405
$ws.Range("J5").Value = 2090.5444  # was 2071.5405
$ws.Range("L5").Value = 2090.5444  # was 2071.5405
$ws.Range("N5").Value = -2314.5444  # was -2295.5405

# Row 132
$ws.Range("H132").Value = 2617.12  # was 2443.3
$ws.Range("I132").Value = 2463.3333  # was 2383.9546
$ws.Range("J132").Value = 3424.5  # was 2606.5
$ws.Range("K132").Value = 7389.999899999999  # was 7151.8638
$ws.Range("L132").Value = 10273.5  # was 7819.5
$ws.Range("M132").Value = -4859.999899999999  # was -4621.8638
$ws.Range("N132").Value = -15333.5  # was -12879.5

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 5006415.5  # was 6005000.5
$ws.Range("I2").Value = 0  # was 1000
$ws.Range("J2").Value = 5006415.5  # was 7506000.5
$ws.Range("K2").Value = 0  # was 1000
$ws.Range("L2").Value = 5006415.5  # was 7506000.5
$ws.Range("M2").ClearContents()  # was -888
$ws.Range("N2").Value = -5006639.5  # was -7506224.5
